# Applies the Book2.xlsx edits described in the commit:
#  - Shift the DEBIT/CREDIT header labels one column to the left
#    (D1 -> C1, E1 -> D1) and leave the trailing E1 cell blank
#    (while keeping its bold header formatting).
#  - Widen column B to fit the "PARTICULARS" header.
#  - Update the sheet view: new zoom level and active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the current header text for the cells we are about to shift.
$debitText  = $ws.Range("D1").Text
$creditText = $ws.Range("E1").Text

# Copy formatting (bold header style) one column to the left before the
# values move, so the destination cells keep the same look.
$ws.Range("D1").Copy($ws.Range("C1")) | Out-Null
$ws.Range("E1").Copy($ws.Range("D1")) | Out-Null

# Now move the actual text: D1 -> C1, E1 -> D1, then blank out E1
# (its formatting was already copied onto D1 above, and E1 itself keeps
# its own existing bold formatting).
$ws.Range("C1").Value = $debitText
$ws.Range("D1").Value = $creditText
$ws.Range("E1").ClearContents() | Out-Null

# Widen column B (PARTICULARS).
$ws.Columns.Item(2).ColumnWidth = 17.88671875

# Update the view: select C4 and set the zoom level.
$ws.Range("C4").Select() | Out-Null
$excel.ActiveWindow.Zoom = 119
